{"js": "// Insert the clause \" v\u00e0 file JSON cho Postman\" right after the existing\n// \"back end \" text (and before \"h\u1ecd \u0111\u00e3 build s\u1eb5n, ...\"), inside the\n// paragraph that starts with \"T\u1ea3i Docker desktop, ...\".\n//\n// Net effect (confirmed against the OOXML diff):\n//   \"...ho\u1ea1t back end h\u1ecd \u0111\u00e3 build s\u1eb5n...\"\n// becomes\n//   \"...ho\u1ea1t back end v\u00e0 file JSON cho Postman h\u1ecd \u0111\u00e3 build s\u1eb5n...\"\n\nconst body = context.document.body;\n\n// \"back end \" (with the trailing space) is its own run in the original\n// document, immediately followed by the run containing \"h\u1ecd\". Searching\n// for it (case sensitive, whole match) gives us a precise anchor to\n// insert after, without disturbing the rest of the paragraph.\nconst results = body.search(\"back end \", { matchCase: true, matchWholeWord: false });\nresults.load(\"items,text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find anchor text \"back end \" in the document.');\n}\n\nconst anchor = results.items[0];\nanchor.insertText(\"v\u00e0 file JSON cho Postman \", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Insert the clause \"v\u00e0 file JSON cho Postman \" right after the existing\n# \"back end \" text (and before \"h\u1ecd \u0111\u00e3 build s\u1eb5n, ...\"), inside the\n# paragraph that starts with \"T\u1ea3i Docker desktop, ...\".\n#\n# Net effect (confirmed against the OOXML diff):\n#   \"...ho\u1ea1t back end h\u1ecd \u0111\u00e3 build s\u1eb5n...\"\n# becomes\n#   \"...ho\u1ea1t back end v\u00e0 file JSON cho Postman h\u1ecd \u0111\u00e3 build s\u1eb5n...\"\n\n$d = $word.ActiveDocument\n\n# \"back end \" (with the trailing space) is its own run in the original\n# document, immediately followed by the run containing \"h\u1ecd\". Find it\n# (case sensitive, exact text) so we can anchor the insertion precisely.\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"back end \"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$found = $find.Execute()\n\nif ($found) {\n    # $rng now spans the found \"back end \" text; InsertAfter adds the new\n    # text immediately following it, preserving the surrounding formatting.\n    $rng.InsertAfter(\"v\u00e0 file JSON cho Postman \")\n} else {\n    throw 'Could not find anchor text \"back end \" in the document.'\n}\n"}
